$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must stay literal text
# (preserve formatting such as trailing zeros, e.g. "1.00", "0.150").
# Pre-format them as Text so Excel does not normalize the value to a number.
$textCells = @("D9", "D13", "D15", "D20", "D21", "D27", "D36", "D43", "D45")
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

$ws.Range("D2").Value = "90.700.03"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.128.19"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "238.14"
$ws.Range("E5").Value = "  +8.36%  "
$ws.Range("D6").Value = "630.46"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "1.06"
$ws.Range("E7").Value = "  +10.83%  "
$ws.Range("D8").Value = "0.357"
$ws.Range("E8").Value = "  -8.29%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.123.46"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "0.725"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "36.70"
$ws.Range("E13").Value = "  +6.01%  "
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").Value = "0.0000243"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "90.309.75"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "3.694.86"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "3.126.93"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "14.50"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("D21").Value = "0.0000212"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "452.77"
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("D23").Value = "5.71"
$ws.Range("E23").Value = "  +10.35%  "
$ws.Range("D24").Value = "9.08"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").Value = "6.07"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "90.71"
$ws.Range("E26").Value = "  +5.25%  "
$ws.Range("D27").Value = "12.60"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("D28").Value = "3.275.95"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "9.85"
$ws.Range("E30").Value = "  +8.62%  "
$ws.Range("D31").Value = "0.161"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").Value = "27.67"
$ws.Range("E32").Value = "  +17.18%  "
$ws.Range("D33").Value = "0.201"
$ws.Range("E33").Value = "  +34.34%  "
$ws.Range("D34").Value = "3.83"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "513.17"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.150"
$ws.Range("E36").Value = "  +5.29%  "
$ws.Range("D37").Value = "7.15"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D41").Value = "0.429"
$ws.Range("E41").Value = "  +13.28%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "0.0869"
$ws.Range("E42").Value = "  +5.40%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "22.20"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D45").Value = "3.40"
$ws.Range("E45").Value = "  +43.74%  "
$ws.Range("D46").Value = "1.96"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "0.701"
$ws.Range("E47").Value = "  +13.32%  "
$ws.Range("D48").Value = "149.32"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "45.62"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").Value = "4.56"
$ws.Range("E50").Value = "  +9.14%  "
$ws.Range("D51").Value = "1.36"
$ws.Range("E51").Value = "  +4.73%  "
